$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jpy")
$ws.Unprotect()
$ws.Range("A30").Value = "hello"
$ws.Protect($false, $true, $true, $true)
